# [FEATURE] Add tablet settings screens
# Adds a new "Tablet Album und Kochbuch Einstellungen" time-tracking entry
# to the Arbeitsmatrix sheet, inserting it (plus two extra blank buffer
# rows) right before the existing blank-row buffer / summary block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# The sheet currently has a blank buffer row at 151/152 immediately above
# the "Stunden insgesamt" summary block (row 153). Insert 3 fresh blank
# rows at 152 so that:
#   - row 151 stays the existing blank buffer row
#   - rows 152-154 are brand new blank rows (formatted like row 151)
#   - the old blank row 152 becomes row 155
#   - the old summary rows 153-156 become rows 156-159
$ws.Rows("152:154").Insert()

# Copy the formatting from the row above (row 150, the last real entry)
# onto the new row 152 so styles/number formats match the other entries.
$ws.Range("A150:G150").Copy()
$ws.Range("A152:G152").PasteSpecial(-4122)
$ws.Range("I150:K150").Copy()
$ws.Range("I152:K152").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new entry's data.
$ws.Range("A152").Value = 22
$ws.Range("B152").Value = "Interface Design"
$ws.Range("C152").Value = "MockUps"
$ws.Range("D152").Value = "[FEATURE]"
$ws.Range("E152").Value = "Tablet Album und Kochbuch Einstellungen"
$ws.Range("F152").Value = 44494
$ws.Range("G152").Value = 44481
$ws.Range("J152").Value = 0.5
$ws.Range("K152").Value = 0.60416666666666663
$ws.Range("I152").Formula = "=ROUNDUP(((SUM(K152-J152)*24*60/60)/0.25),0)*0.25"

# Keep the sheet's selection pointing at the new bottom-right corner of
# the used range (the summary block now ends on row 159).
$ws.Range("A1:I159").Select()

$wb.Save()
